# Thank You slide (slide 12): simplify messaging to focus on
# "Questions & Discussion" / team-communication / thank-you copy,
# per commit "fix(presentation): improve Thank You and Q&A slide readability".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)          # "Subtitle 2" placeholder (subTitle)
$tr = $sh.TextFrame.TextRange

# Paragraphs (1-indexed) before the edit:
#   1: "Ready to Improve Team Communication?"
#   2: (blank)
#   3: "💬 Start Phase 1 This Week"
#   4: "📈 Measure the Impact"
#   5: "🚀 Prepare for Phase 2"
#   6: (blank)
#   7: "Let's build better software, together."
#
# Drop paragraphs 4 and 5 ("Measure the Impact" / "Prepare for Phase 2"),
# deleting from the back so earlier indices stay valid.
$tr.Paragraphs(5, 1).Delete()
$tr.Paragraphs(4, 1).Delete()

# Remaining layout now:
#   1: "Ready to Improve Team Communication?"
#   2: (blank)
#   3: "💬 Start Phase 1 This Week"
#   4: (blank)
#   5: "Let's build better software, together."

# Replace paragraph 3's text with the new single message.
$tr.Paragraphs(3, 1).Text = "Let's improve our team communication together."

# Title line -> "Questions & Discussion", bumped up to 20pt.
$tr.Paragraphs(1, 1).Text = "Questions & Discussion"
$tr.Paragraphs(1, 1).Font.Size = 20

# Closing line -> "Thank you for your attention!"
$tr.Paragraphs(5, 1).Text = "Thank you for your attention!"
